$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") for data rows 2-10: 45243 -> 45244
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 45244
}
